$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "CAN" -> "CNA" in cell A7
$ws.Range("A7").Value = "CNA"

# Update selection to B10 (active cell), even though it has no data
$ws.Range("B10").Select()
